$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old 2-column lookup table (A1:B5) before laying out the new
# wide table.
$ws.Range("A1:B5").ClearContents()

# Header row (row 1): first write all eleven header labels...
$headers = @(
    "MIGRATION DATE",
    "FINANCIAL INSTITUTION NAME",
    "ENTITY ID",
    "ADDRESS",
    "CITY",
    "STATE",
    "ZIP CODE",
    "PHONE #",
    "PROJECT COORDINATOR",
    "CERTIFICATION REQUIRED (Yes or No)",
    "CERTIFICATION COORDINATOR"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ...then reuse the existing bold / centered / bordered header style (already
# applied to A1:B1) across the rest of the header row C1:K1.
$ws.Cells.Item(1, 1).Copy()
$ws.Range("C1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows: one value per row, staggered across the matching column.
# Force text format on cells whose content would otherwise be
# auto-coerced to a date/number by Excel (the migration date, zip code and
# phone number all need to stay literal text).
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-10-16"

$ws.Cells.Item(3, 2).Value = "YYY"

$ws.Cells.Item(4, 3).Value = "123ABX007"

$ws.Cells.Item(5, 4).Value = "Karapakkam"

$ws.Cells.Item(6, 5).Value = "Chennai"

$ws.Cells.Item(7, 6).Value = "Tamil Nadu"

$ws.Cells.Item(8, 7).NumberFormat = "@"
$ws.Cells.Item(8, 7).Value = "600117"

$ws.Cells.Item(9, 8).NumberFormat = "@"
$ws.Cells.Item(9, 8).Value = "9911991100"

$ws.Cells.Item(10, 9).Value = "Sam"

$ws.Cells.Item(11, 10).Value = "Yes"

# Touch the bottom-right corner of the table so the sheet's used range /
# dimension extends through row 12, column K (K12 itself stays blank).
$ws.Cells.Item(12, 11).NumberFormat = "General"
